$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": update PORCELANATO value for the client row ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M22").Value = -2156.54

# --- Sheet "VENTA MENSUAL": update agosto (August) values ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F22").Value = -2156.54
$ws2.Range("F23").Value = 11812.64
$ws2.Columns.Item(6).ColumnWidth = 13.166666666666666

# --- Sheet "CUMPLIMIENTO MENSUAL": update PORCELANATO row and TOTAL row ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D16").Value = 9323.23
$ws3.Range("E16").Value = 29453.24
$ws3.Range("F16").Value = 0.2404352433318453
$ws3.Range("D19").Value = 13527.95
$ws3.Range("E19").Value = 45860.27762291768
$ws3.Range("F19").Value = 0.2277884109607544
$ws3.Columns.Item(4).ColumnWidth = 12.166666666666666
